# Apply updated K (strikeout) values to column G, regenerated from source data
# (replacing the old "Strike#" values) for data/save_data/2021/chafin_andrew.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 0
    3 = 2
    4 = 1
    5 = 0
    6 = 1
    7 = 0
    8 = 2
    9 = 1
    10 = 0
    11 = 0
    12 = 1
    13 = 4
    14 = 2
    15 = 2
    16 = 1
    17 = 0
    18 = 2
    19 = 0
    20 = 1
    21 = 2
    22 = 0
    23 = 1
    24 = 2
    25 = 2
    26 = 1
    27 = 0
    28 = 1
    29 = 1
    30 = 0
    31 = 1
    32 = 1
    33 = 1
    34 = 2
    35 = 1
    36 = 1
    37 = 0
    38 = 1
    39 = 0
    40 = 0
    41 = 1
    42 = 1
    43 = 1
    44 = 2
    45 = 0
    46 = 0
    47 = 1
    48 = 1
    49 = 1
    50 = 1
    51 = 0
    52 = 0
    53 = 0
    54 = 1
    55 = 2
    56 = 1
    57 = 1
    58 = 1
    59 = 1
    60 = 1
    61 = 0
    62 = 1
    63 = 0
    64 = 2
    65 = 0
    66 = 1
    67 = 1
    68 = 0
    69 = 2
    70 = 4
    71 = 2
    72 = 0
    74 = 2
    77 = 1
    78 = 1
    79 = 1
    80 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $kValues[$row]
}

Write-Host "Updated" $kValues.Count "K values in column G"
